# Update the "想去人数" (want-to-go count) column F on several sheets.
# Each listed cell's numeric value is incremented by 1 (matching a refreshed
# snapshot of the live counts that feed this static export).

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        "F5"  = 404
        "F6"  = 791
        "F8"  = 1129
        "F9"  = 309
        "F14" = 505
        "F15" = 140
        "F18" = 2870
        "F27" = 587
        "F29" = 18
        "F30" = 56
        "F31" = 284
        "F32" = 1061
    }
    "本地生活" = @{
        "F5" = 2424
    }
    "全部类型" = @{
        "F3"  = 2424
        "F11" = 404
        "F12" = 791
        "F15" = 1129
        "F16" = 309
        "F20" = 505
        "F23" = 2870
        "F36" = 18
        "F38" = 284
        "F44" = 1061
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $updates[$sheetName]
    foreach ($addr in $cellMap.Keys) {
        $ws.Range($addr).Value = $cellMap[$addr]
    }
}
